# Weekly update: insert a new price record for the most recent week
# (Feria Lagunitas de Puerto Montt - Espárragos) and push the existing
# historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 14; rows 14..46 shift down to 15..47
# and inherit the formatting (e.g. the date style on column D).
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with this week's record.
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C14").Value = "Los Lagos"
$ws.Range("D14").Value = 44838
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 300000000
$ws.Range("G14").Value = "Espárragos"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 300
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 2000
$ws.Range("N14").Value = "`$/kilo"
$ws.Range("O14").Value = "Provincia de Linares"
$ws.Range("P14").Value = 2000
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
